# Weekly price update: insert one new price record for
# "Comercializadora del Agro de Limarí - Ají" as a new row 431,
# pushing the existing rows 431-460 down to 432-461.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 431 (shifts old rows 431:460 -> 432:461,
# and grows the sheet's used range from R460 to R461 automatically).
$ws.Rows("431:431").Insert()

# Columns that are constant for every record in this sheet - copy them
# from the row immediately below (the row that used to be 431).
$ws.Cells.Item(431, 1).Value = $ws.Cells.Item(432, 1).Value2    # A  Mercado ID
$ws.Cells.Item(431, 2).Value = $ws.Cells.Item(432, 2).Value2    # B  Mercado
$ws.Cells.Item(431, 3).Value = $ws.Cells.Item(432, 3).Value2    # C  Region
$ws.Cells.Item(431, 5).Value = $ws.Cells.Item(432, 5).Value2    # E  Codreg
$ws.Cells.Item(431, 6).Value = $ws.Cells.Item(432, 6).Value2    # F  Categoria ID
$ws.Cells.Item(431, 7).Value = $ws.Cells.Item(432, 7).Value2    # G  Categoria
$ws.Cells.Item(431, 15).Value = $ws.Cells.Item(432, 15).Value2  # O  Origen
$ws.Cells.Item(431, 18).Value = $ws.Cells.Item(432, 18).Value2  # R  Clasificacion

# New record's own values.
$ws.Cells.Item(431, 4).Value = 45147               # D  Fecha
$ws.Cells.Item(431, 8).Value = "Americana (o)"     # H  Variedad
$ws.Cells.Item(431, 9).Value = "Primera"           # I  Calidad
$ws.Cells.Item(431, 10).Value = 200                # J  Volumen
$ws.Cells.Item(431, 11).Value = 29000              # K  Precio minimo
$ws.Cells.Item(431, 12).Value = 30000              # L  Precio maximo
$ws.Cells.Item(431, 13).Value = 29500              # M  Precio promedio ponderado
$ws.Cells.Item(431, 14).Value = "$/caja 25 kilos"  # N  Unidad de comercializacion
$ws.Cells.Item(431, 16).Value = 1180               # P  Precio $/Kg
$ws.Cells.Item(431, 17).Value = 25                 # Q  Kg o Unidades
